$wb = $excel.ActiveWorkbook

# --- resources sheet: clear scrolled view, keep selection B3 ---
$resources = $wb.Worksheets.Item("resources")
$resources.Activate()
$resources.Range("A1").Select()
$resources.Range("B3").Select()

# --- sps sheet: fix the combo_pairs text for the "publisher" row ---
$sps = $wb.Worksheets.Item("sps")
$sps.Range("Q6").Value = "context-type-value,version"

# --- ops sheet: add a "definition" column with a hyperlink to the operation definition ---
$ops = $wb.Worksheets.Item("ops")
$ops.Columns.Item(2).Insert()
$ops.Range("B1").Value = "definition"
$ops.Hyperlinks.Add($ops.Range("B2"), "http://fhir.org/guides/argonaut-questionnaire/OperationDefinition/next-question")
$ops.Columns.Item(2).ColumnWidth = 83

# --- final view state: "ops" sheet active, selection on B9; "sps" selection on Q8 ---
$sps.Activate()
$sps.Range("Q8").Select()
$ops.Activate()
$ops.Range("B9").Select()
